$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Estado de completado" comment for the "Lectura de base de datos" row
$ws.Range("E4").Value = "Voy a hacer el la lectura de la entidad menu"

# Resize column E (best fit) so the new, longer text is fully visible
$ws.Columns.Item(5).ColumnWidth = 39.25

# Move the active cell/selection
$ws.Range("F5").Select() | Out-Null
